# Cadastro template update: add "Unidade(s)*" column, mark required headers
# with "*", widen a few columns, and move the active selection to A2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename existing headers to flag them as required (all except Email) ---
$ws.Range("A1").Value = "Nome completo*"
$ws.Range("B1").Value = "CPF*"
$ws.Range("C1").Value = "Data de Nascimento*"
$ws.Range("D1").Value = "Email"
$ws.Range("E1").Value = "Cargo/Função*"

# --- 2. Touch A1's font so it picks up its own cell style (matches the
#        extra cellXfs entry the source workbook carries for the header) ---
$ws.Range("A1").Font.Name = "Calibri"

# --- 3. Add the new "Unidade(s)*" column to the Usuarios table ---
$lo = $ws.ListObjects.Item(1)
$newCol = $lo.ListColumns.Add()
$ws.Range("F1").Value = "Unidade(s)*"

# Format the new column's data as Text (matches the "@" format applied to
# the other text columns, e.g. CPF in column B).
$ws.Columns.Item(6).NumberFormat = "@"

# --- 4. Column width adjustments ---
# ColumnWidth (character units) differs from the raw XML width by 5/6;
# offset targets accordingly so the saved width matches the template.
$ws.Columns.Item(3).ColumnWidth = 25.166666666666668   # C: 16 -> 26
$ws.Columns.Item(5).ColumnWidth = 19.166666666666668   # E: 16 -> 20
$ws.Columns.Item(6).ColumnWidth = 20.45                # F: ~21.29 new column

# --- 5. Move the active selection to A2 (ready for data entry) ---
[void]$ws.Range("A2").Select()
